$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MO")

# Row 15 - Gross Margin
$ws.Range("D15").Value = 0.7485
$ws.Range("E15").Value = 0.6483
$ws.Range("F15").Value = 0.5656
$ws.Range("G15").Value = 0.5062

# Row 16 - EBIT Margin
$ws.Range("D16").Value = 0.4919
$ws.Range("E16").Value = 0.4572
$ws.Range("F16").Value = 0.425
$ws.Range("G16").Value = 0.4112

# Row 17 - EBT margin
$ws.Range("D17").Value = 0.1434
$ws.Range("E17").Value = 0.0576
$ws.Range("F17").Value = 0.0554
$ws.Range("G17").Value = 0.0305

# Row 18 - Net Profit Margin
$ws.Range("D18").Value = 0.0332
$ws.Range("E18").Value = -0.0402
$ws.Range("F18").Value = -0.0354
$ws.Range("G18").Value = -0.0518

# Row 19 - Free Cash Flow Margin
$ws.Range("B19").Value = 0.3913
$ws.Range("D19").Value = 0.3744
$ws.Range("E19").Value = 0.4405
$ws.Range("F19").Value = 0.3432
$ws.Range("G19").Value = 0.3023

# Row 27 - EBITDA Margin
$ws.Range("D27").Value = 0.6606
$ws.Range("E27").Value = 0.5634
$ws.Range("F27").Value = 0.4804
$ws.Range("G27").Value = 0.4202

# Row 28 - Operating Cash Flow Margin
$ws.Range("D28").Value = 0.3858
$ws.Range("E28").Value = 0.4524
$ws.Range("F28").Value = 0.3538
$ws.Range("G28").Value = 0.3121
